$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data for I2:J61 (60 rows x 2 columns)
$data = New-Object 'object[,]' 60,2
$data[0,0] = 9; $data[0,1] = 9
$data[1,0] = 8; $data[1,1] = 8
$data[2,0] = 9; $data[2,1] = 9
$data[3,0] = 8; $data[3,1] = 8
$data[4,0] = 8; $data[4,1] = 8
$data[5,0] = 7; $data[5,1] = 8
$data[6,0] = 7; $data[6,1] = 7
$data[7,0] = 6; $data[7,1] = 7
$data[8,0] = 7; $data[8,1] = 7
$data[9,0] = 7; $data[9,1] = 7
$data[10,0] = 7; $data[10,1] = 7
$data[11,0] = 8; $data[11,1] = 8
$data[12,0] = 7; $data[12,1] = 8
$data[13,0] = 10; $data[13,1] = 10
$data[14,0] = 4; $data[14,1] = 5
$data[15,0] = 5; $data[15,1] = 6
$data[16,0] = 8; $data[16,1] = 9
$data[17,0] = 8; $data[17,1] = 8
$data[18,0] = 9; $data[18,1] = 9
$data[19,0] = 7; $data[19,1] = 7
$data[20,0] = 8; $data[20,1] = 8
$data[21,0] = 8; $data[21,1] = 9
$data[22,0] = 4; $data[22,1] = 5
$data[23,0] = 4; $data[23,1] = 5
$data[24,0] = 6; $data[24,1] = 6
$data[25,0] = 11; $data[25,1] = 11
$data[26,0] = 9; $data[26,1] = 9
$data[27,0] = 8; $data[27,1] = 8
$data[28,0] = 9; $data[28,1] = 9
$data[29,0] = 7; $data[29,1] = 7
$data[30,0] = 9; $data[30,1] = 9
$data[31,0] = 6; $data[31,1] = 6
$data[32,0] = 9; $data[32,1] = 9
$data[33,0] = 6; $data[33,1] = 6
$data[34,0] = 6; $data[34,1] = 7
$data[35,0] = 5; $data[35,1] = 6
$data[36,0] = 7; $data[36,1] = 7
$data[37,0] = 5; $data[37,1] = 5
$data[38,0] = 7; $data[38,1] = 7
$data[39,0] = 7; $data[39,1] = 7
$data[40,0] = 8; $data[40,1] = 8
$data[41,0] = 1; $data[41,1] = 1
$data[42,0] = 5; $data[42,1] = 5
$data[43,0] = 8; $data[43,1] = 8
$data[44,0] = 5; $data[44,1] = 6
$data[45,0] = 5; $data[45,1] = 5
$data[46,0] = 7; $data[46,1] = 7
$data[47,0] = 8; $data[47,1] = 8
$data[48,0] = 8; $data[48,1] = 8
$data[49,0] = 8; $data[49,1] = 8
$data[50,0] = 7; $data[50,1] = 7
$data[51,0] = 7; $data[51,1] = 7
$data[52,0] = 7; $data[52,1] = 7
$data[53,0] = 8; $data[53,1] = 8
$data[54,0] = 8; $data[54,1] = 8
$data[55,0] = 3; $data[55,1] = 3
$data[56,0] = 4; $data[56,1] = 4
$data[57,0] = 6; $data[57,1] = 6
$data[58,0] = 6; $data[58,1] = 6
$data[59,0] = 4; $data[59,1] = 4

$ws.Range("I2:J61").Value = $data
